$d = $word.ActiveDocument

# Merge the split runs of the title paragraph into a single run.
$d.Content.Find.Execute("Sigma Notation: Answers", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Sigma Notation: Answers", 2)

# Merge the split runs of the author paragraph into a single run.
$d.Content.Find.Execute("Ifan Howells-Baines, Mark Toner", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ifan Howells-Baines, Mark Toner", 2)
